# Updated cryptos list on Sat May 27 22:13:18 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row (rows 2-51) on the active sheet to the latest scraped snapshot.
#
# A few Price strings look like plain decimals (e.g. "1.020", "66.60") and
# Excel's COM type-sniffing would silently coerce them to numbers (dropping
# the significant trailing zero). To keep them as literal text -- matching
# the source data -- we briefly force Text format, assign the value, then
# restore the Normal cell style so no visible formatting changes remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "27.264.16"
$ws.Range("E2").Value = "  +1.48%  "

$ws.Range("D3").Value = "1.861.70"
$ws.Range("E3").Value = "  +1.13%  "

Set-TextValue $ws.Range("D4") "1.020"
$ws.Range("E4").Value = "  +1.32%  "

$ws.Range("E5").Value = "  +1.10%  "

Set-TextValue $ws.Range("D6") "1.019"
$ws.Range("E6").Value = "  +1.40%  "

Set-TextValue $ws.Range("D7") "0.4797"
$ws.Range("E7").Value = "  +1.97%  "

Set-TextValue $ws.Range("D8") "0.3726"
$ws.Range("E8").Value = "  +1.79%  "

Set-TextValue $ws.Range("D9") "0.07326"
$ws.Range("E9").Value = "  +2.55%  "

Set-TextValue $ws.Range("D10") "0.9369"
$ws.Range("E10").Value = "  +1.53%  "

Set-TextValue $ws.Range("D11") "20.34"
$ws.Range("E11").Value = "  +4.09%  "

Set-TextValue $ws.Range("D12") "0.07871"
$ws.Range("E12").Value = "  +2.70%  "

$ws.Range("D13").Value = "1.892.66"
$ws.Range("E13").Value = "  +2.35%  "

Set-TextValue $ws.Range("D14") "5.425"
$ws.Range("E14").Value = "  +2.62%  "

Set-TextValue $ws.Range("D15") "6.531"
$ws.Range("E15").Value = "  +2.07%  "

Set-TextValue $ws.Range("D16") "90.09"
$ws.Range("E16").Value = "  +2.06%  "

Set-TextValue $ws.Range("D17") "1.021"
$ws.Range("E17").Value = "  +1.32%  "

Set-TextValue $ws.Range("D18") "0.000008747"
$ws.Range("E18").Value = "  +1.35%  "

$ws.Range("E19").Value = "  +1.35%  "

Set-TextValue $ws.Range("D20") "14.78"
$ws.Range("E20").Value = "  +2.26%  "

$ws.Range("D21").Value = "27.301.86"
$ws.Range("E21").Value = "  +1.49%  "

Set-TextValue $ws.Range("D22") "5.113"
$ws.Range("E22").Value = "  +2.08%  "

$ws.Range("E23").Value = "  +0.68%  "

Set-TextValue $ws.Range("D24") "1.950"
$ws.Range("E24").Value = "  +1.33%  "

Set-TextValue $ws.Range("D25") "153.88"
$ws.Range("E25").Value = "  +1.47%  "

$ws.Range("E26").Value = "  +1.55%  "

Set-TextValue $ws.Range("D27") "2.001"
$ws.Range("E27").Value = "  -0.44%  "

Set-TextValue $ws.Range("D28") "115.72"
$ws.Range("E28").Value = "  +1.40%  "

Set-TextValue $ws.Range("D29") "4.989"
$ws.Range("E29").Value = "  +2.29%  "

Set-TextValue $ws.Range("D30") "0.08895"
$ws.Range("E30").Value = "  +0.88%  "

Set-TextValue $ws.Range("D31") "3.348"
$ws.Range("E31").Value = "  +4.29%  "

Set-TextValue $ws.Range("D32") "1.191"
$ws.Range("E32").Value = "  +1.26%  "

Set-TextValue $ws.Range("D33") "4.567"
$ws.Range("E33").Value = "  +1.97%  "

Set-TextValue $ws.Range("D34") "0.7401"
$ws.Range("E34").Value = "  -0.68%  "

Set-TextValue $ws.Range("D35") "2.679"
$ws.Range("E35").Value = "  -2.83%  "

Set-TextValue $ws.Range("D36") "1.124"
$ws.Range("E36").Value = "  +3.50%  "

Set-TextValue $ws.Range("D37") "0.02034"
$ws.Range("E37").Value = "  +4.88%  "

Set-TextValue $ws.Range("D38") "0.05261"
$ws.Range("E38").Value = "  +1.02%  "

Set-TextValue $ws.Range("D39") "0.5361"

Set-TextValue $ws.Range("D40") "7.127"
$ws.Range("E40").Value = "  +2.36%  "

Set-TextValue $ws.Range("D41") "0.1533"
$ws.Range("E41").Value = "  +1.67%  "

Set-TextValue $ws.Range("D42") "8.349"
$ws.Range("E42").Value = "  +2.43%  "

$ws.Range("E43").Value = "  +1.71%  "

Set-TextValue $ws.Range("D44") "0.4803"
$ws.Range("E44").Value = "  +2.29%  "

Set-TextValue $ws.Range("D45") "1.020"
$ws.Range("E45").Value = "  +1.50%  "

Set-TextValue $ws.Range("D46") "102.85"
$ws.Range("E46").Value = "  +1.07%  "

Set-TextValue $ws.Range("D47") "1.636"
$ws.Range("E47").Value = "  +2.43%  "

Set-TextValue $ws.Range("D48") "66.60"
$ws.Range("E48").Value = "  +0.97%  "

Set-TextValue $ws.Range("D49") "0.06077"
$ws.Range("E49").Value = "  +0.57%  "

Set-TextValue $ws.Range("D50") "0.9011"
$ws.Range("E50").Value = "  +1.40%  "

Set-TextValue $ws.Range("D51") "36.74"
$ws.Range("E51").Value = "  +1.53%  "
